# feat: Add SFX of MainMenu resource file
#
# Adds a new "SFX" worksheet (right after the existing "BGM" sheet) that
# lists a single sound-effect resource row, mirroring the layout/headers
# of the BGM sheet and its "BGM_..." CONCAT-based FileName formula.

$wb  = $excel.ActiveWorkbook
$bgm = $wb.Worksheets.Item("BGM")

# Insert the new sheet immediately after "BGM" and name it "SFX".
$sfx = $wb.Worksheets.Add([System.Type]::Missing, $bgm)
$sfx.Name = "SFX"

# Header row - identical column headers to the BGM sheet.
$sfx.Cells.Item(1,1).Value = "Team"
$sfx.Cells.Item(1,2).Value = "Class"
$sfx.Cells.Item(1,3).Value = "Situation"
$sfx.Cells.Item(1,4).Value = "Source"
$sfx.Cells.Item(1,5).Value = "Single / Multiple"
$sfx.Cells.Item(1,6).Value = "Audio Normalization"
$sfx.Cells.Item(1,7).Value = "FileName"

# Row 2 - the new MainMenu / TitleScreen / MenuSelect SFX entry.
$sfx.Cells.Item(2,1).Value = "MainMenu"
$sfx.Cells.Item(2,2).Value = "TitleScreen"
$sfx.Cells.Item(2,3).Value = "MenuSelect"
$sfx.Cells.Item(2,4).Value = "BlipLow.wav"
$sfx.Cells.Item(2,5).Value = "S"
$sfx.Cells.Item(2,6).Value = "O"
$sfx.Cells.Item(2,7).Formula = '=CONCAT("SFX_",A2,"_",B2,"_",C2,"_",D2)'

# Trailing styled-but-empty cell (row 17), same "blank row" font style
# used for the blank row at the bottom of the BGM sheet.
$sfx.Cells.Item(17,3).Font.Color = 0

# Column widths close to the source sheet's autofit widths.
$sfx.Columns.Item(1).ColumnWidth = 24.285714
$sfx.Columns.Item(2).ColumnWidth = 24.285714
$sfx.Columns.Item(3).ColumnWidth = 28.856028328125024
$sfx.Columns.Item(4).ColumnWidth = 24.1411826921875
$sfx.Columns.Item(5).ColumnWidth = 14.141178692187443
$sfx.Columns.Item(6).ColumnWidth = 17.570867564062464
$sfx.Columns.Item(7).ColumnWidth = 74.28573400000026

# Zoom level used on the new sheet.
$sfx.Application.ActiveWindow.Zoom = 154
